$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Rana Abo-Zaid, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G4").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Range("G5").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Nada Gouda'
$ws.Range("G6").Value = 'Dr. Nada Mohammad, Dr. Kerelos Zareef'
$ws.Range("G8").Value = 'Dr. Aya Saeed, Dr. Amal Awwad'
$ws.Range("G9").Value = 'Dr. Aya Essam, Dr. Nourhan Mohammad'
$ws.Range("G10").Value = 'Dr. Marina Youhanna, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Amany Raafat'
$ws.Range("G11").Value = 'Dr. Sarah Mahdy, Dr. Alaa Ashraf'
$ws.Range("G18").Value = 'Dr. Remon, Dr. Aya Hanafy, Dr. Yasmin, Dr. Shorok Mohammad'
$ws.Range("G19").Value = 'Dr. Naema Gomaa, Dr. Salma Hassan, Dr. Monica, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Yassmen Ahmad, Dr. Remon, Dr. Nardine, Dr. Wafaa Ebida'
$ws.Range("G20").Value = 'Dr. Youstina Magdy, Dr. Aya Hanafy, Dr. Nardine, Dr. Yassmen Ahmad, Dr. Remon, Dr. Marina Sorial, Dr. Wafaa Ebida'
$ws.Range("G21").Value = 'Dr. Monica, Dr. Neveen Nashaat, Dr. Shorok Mohammad, Dr. Yassmen Ahmad, Dr. Yasmin'
$ws.Range("G22").Value = 'Dr. Remon, Dr. Monica, Dr. Naema Gomaa, Dr. Wafaa Ebida'
$ws.Range("G24").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Rana Abo-Zaid, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G25").Value = 'Dr. Alshimaa Atef, Dr. Manar Montaser, Administrator, Dr. Gehan Adel'
$ws.Range("G26").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Range("G27").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Nada Gouda'
$ws.Range("G28").Value = 'Dr. Nada Mohammad, Dr. Kerelos Zareef'
$ws.Range("G30").Value = 'Dr. Aya Saeed, Dr. Amal Awwad'
$ws.Range("G32").Value = 'Dr. Marina Youhanna, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Amany Raafat'
$ws.Range("G33").Value = 'Dr. Sarah Mahdy, Dr. Alaa Ashraf'
$ws.Range("G40").Value = 'Dr. Remon, Dr. Aya Hanafy, Dr. Yasmin, Dr. Shorok Mohammad'
$ws.Range("G41").Value = 'Dr. Naema Gomaa, Dr. Salma Hassan, Dr. Monica, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Yassmen Ahmad, Dr. Remon, Dr. Nardine, Dr. Wafaa Ebida'
$ws.Range("G42").Value = 'Dr. Youstina Magdy, Dr. Aya Hanafy, Dr. Nardine, Dr. Yassmen Ahmad, Dr. Remon, Dr. Marina Sorial, Dr. Wafaa Ebida'
$ws.Range("G43").Value = 'Dr. Monica, Dr. Neveen Nashaat, Dr. Shorok Mohammad, Dr. Yassmen Ahmad, Dr. Yasmin'
$ws.Range("G44").Value = 'Dr. Remon, Dr. Monica, Dr. Naema Gomaa, Dr. Wafaa Ebida'
$ws.Range("G46").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Nahla Nagiub, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud'
$ws.Range("G48").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad'
$ws.Range("G52").Value = 'Dr. Mariam Nour El-Din, Dr. Shimaa Ashraf'
$ws.Range("G54").Value = 'Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Marwa Mustafa, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya'
$ws.Range("G58").Value = 'Dr. Afaf Abdallah, Dr. Amr Saeed'
$ws.Range("G62").Value = 'Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Shorok Mohammad, Dr. Wafaa Ebida'
$ws.Range("G63").Value = 'Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Wafaa Ebida'
$ws.Range("G64").Value = 'Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Wafaa Ebida'
$ws.Range("G65").Value = 'Dr. Ola Abd Al-Fattah, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Shorok Mohammad, Dr. Remon, Dr. Nardine'
$ws.Range("G66").Value = 'Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Monica, Dr. Maryam Ashraf, Dr. Marina Sorial'
$ws.Range("G68").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Nahla Nagiub, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud'
$ws.Range("G70").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad'
$ws.Range("G74").Value = 'Dr. Mariam Nour El-Din, Dr. Shimaa Ashraf'
$ws.Range("G75").Value = 'Dr. Aya Essam, Dr. Nourhan Mohammad'
$ws.Range("G76").Value = 'Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Marwa Mustafa, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya'
$ws.Range("G80").Value = 'Dr. Afaf Abdallah, Dr. Amr Saeed'
$ws.Range("G84").Value = 'Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Shorok Mohammad, Dr. Wafaa Ebida'
$ws.Range("G85").Value = 'Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Wafaa Ebida'
$ws.Range("G86").Value = 'Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Wafaa Ebida'
$ws.Range("G87").Value = 'Dr. Ola Abd Al-Fattah, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Shorok Mohammad, Dr. Remon, Dr. Nardine'
$ws.Range("G88").Value = 'Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Monica, Dr. Maryam Ashraf, Dr. Marina Sorial'
$ws.Range("G92").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad'
$ws.Range("G93").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Abeer Ragab'
$ws.Range("G96").Value = 'Dr. Sara Nabil, Dr. Mariam Nour El-Din, Dr. Amal Awwad, Dr. Nourhan Mohammad'
$ws.Range("G97").Value = 'Dr. Aya Essam, Dr. Nourhan Mohammad'
$ws.Range("G98").Value = 'Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Marwa Mustafa, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya'
$ws.Range("G106").Value = 'Dr. Youstina Magdy, Dr. Monica, Dr. Neveen Nashaat, Dr. Remon, Dr. Nardine, Dr. Wafaa Ebida'
$ws.Range("G107").Value = 'Dr. Aya Hanafy, Dr. Monica, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Yassmen Ahmad, Dr. Wafaa Ebida'
$ws.Range("G108").Value = 'Dr. Youstina Magdy, Dr. Aya Hanafy, Dr. Nardine, Dr. Yassmen Ahmad, Dr. Remon, Dr. Marina Sorial, Dr. Wafaa Ebida'
$ws.Range("G110").Value = 'Dr. Yassmen Ahmad, Dr. Monica, Dr. Wafaa Ebida'
$ws.Range("G111").Value = 'Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Nourham Mostafa, Dr. Marina Atef, Dr. Monica, Dr. Eman Samir Gabry, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Range("G114").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad'
$ws.Range("G115").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Abeer Ragab'
$ws.Range("G118").Value = 'Dr. Sara Nabil, Dr. Mariam Nour El-Din, Dr. Amal Awwad, Dr. Nourhan Mohammad'
$ws.Range("G119").Value = 'Dr. Aya Essam, Dr. Nourhan Mohammad'
$ws.Range("G120").Value = 'Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Marwa Mustafa, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya'
$ws.Range("G128").Value = 'Dr. Youstina Magdy, Dr. Monica, Dr. Neveen Nashaat, Dr. Remon, Dr. Nardine, Dr. Wafaa Ebida'
$ws.Range("G129").Value = 'Dr. Aya Hanafy, Dr. Monica, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Yassmen Ahmad, Dr. Wafaa Ebida'
$ws.Range("G130").Value = 'Dr. Youstina Magdy, Dr. Aya Hanafy, Dr. Nardine, Dr. Yassmen Ahmad, Dr. Remon, Dr. Marina Sorial, Dr. Wafaa Ebida'
$ws.Range("G131").Value = 'Dr. Marina Atef, Dr. Nardine'
$ws.Range("G132").Value = 'Dr. Yassmen Ahmad, Dr. Monica, Dr. Wafaa Ebida'
$ws.Range("G133").Value = 'Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Nourham Mostafa, Dr. Marina Atef, Dr. Monica, Dr. Eman Samir Gabry, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Range("G134").Value = 'Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Majorelle Magdy, Dr. Amira Sobhy'
$ws.Range("G137").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Nada Gouda'
$ws.Range("G140").Value = 'Dr. Aya Saeed, Dr. Amal Awwad'
$ws.Range("G141").Value = 'Dr. Aya Essam, Dr. Nourhan Mohammad'
$ws.Range("G142").Value = 'Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Merna Said, Dr. Basma Hamed, Dr. Marwa Mustafa'
$ws.Range("G147").Value = 'Dr. Nourham Mostafa, Dr. Nancy Abd Al-Shafy'
$ws.Range("G150").Value = 'Dr. Naema Gomaa, Dr. Salma Hassan, Dr. Monica, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Yassmen Ahmad, Dr. Remon, Dr. Nardine, Dr. Wafaa Ebida'
$ws.Range("G151").Value = 'Dr. Yassmen Ahmad, Dr. Marina Atef, Dr. Monica, Dr. Wafaa Ebida'
$ws.Range("G152").Value = 'Dr. Marina Atef, Dr. Wafaa Ebida'
$ws.Range("G153").Value = 'Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Monica, Dr. Maryam Ashraf, Dr. Marina Sorial'
$ws.Range("G154").Value = 'Dr. Remon, Dr. Naema Gomaa, Dr. Salma Hassan, Dr. Wafaa Ebida'
$ws.Range("G155").Value = 'Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Nourham Mostafa, Dr. Marina Atef, Dr. Monica, Dr. Eman Samir Gabry, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Range("G156").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany, Dr. Alshimaa Atef'
$ws.Range("G159").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Nada Gouda'
$ws.Range("G162").Value = 'Dr. Aya Saeed, Dr. Amal Awwad'
$ws.Range("G163").Value = 'Dr. Aya Essam, Dr. Nourhan Mohammad'
$ws.Range("G164").Value = 'Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Merna Said, Dr. Basma Hamed, Dr. Marwa Mustafa'
$ws.Range("G169").Value = 'Dr. Nourham Mostafa, Dr. Nancy Abd Al-Shafy'
$ws.Range("G172").Value = 'Dr. Naema Gomaa, Dr. Salma Hassan, Dr. Monica, Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Yassmen Ahmad, Dr. Remon, Dr. Nardine, Dr. Wafaa Ebida'
$ws.Range("G173").Value = 'Dr. Yassmen Ahmad, Dr. Marina Atef, Dr. Monica, Dr. Wafaa Ebida'
$ws.Range("G174").Value = 'Dr. Marina Atef, Dr. Wafaa Ebida'
$ws.Range("G175").Value = 'Dr. Eman Mohammad Al, Dr. Aya Hanafy, Dr. Monica, Dr. Maryam Ashraf, Dr. Marina Sorial'
$ws.Range("G176").Value = 'Dr. Remon, Dr. Naema Gomaa, Dr. Salma Hassan, Dr. Wafaa Ebida'
$ws.Range("G177").Value = 'Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Nourham Mostafa, Dr. Marina Atef, Dr. Monica, Dr. Eman Samir Gabry, Dr. Yasmin, Dr. Wafaa Ebida'
